$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF), matching the formatting
# already used by the other header cells (bold, centered, bordered style).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for the new I0 / IF columns, one pair of values per data row (rows 2-36).
$i0 = @(7,3,8,8,7,7,7,8,6,7,8,5,8,6,6,6,9,8,4,6,6,5,8,6,8,8,8,7,7,8,4,4,7,4,3)
$if = @(7,3,9,8,7,7,7,8,6,7,8,5,9,7,6,6,9,8,5,6,6,6,8,7,9,8,8,7,7,8,4,5,7,4,3)

for ($r = 0; $r -lt $i0.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
